# Inserts one new weekly price record for "Terminal Hortofrutícola Agro
# Chillán - Zanahoria" above the existing row 177, shifting the following
# rows (old 177:209) down to (178:210), and fills in the new row's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 177:209 down to 178:210, creating a blank (but formatted) row 177.
$ws.Rows("177:177").Insert()

# Populate the new row 177 with the new weekly record.
$ws.Range("A177").Value = 7
$ws.Range("B177").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C177").Value = "Ñuble"
$ws.Range("D177").Value = 44522
$ws.Range("E177").Value = 16
$ws.Range("F177").Value = 100114013
$ws.Range("G177").Value = "Zanahoria"
$ws.Range("H177").Value = "Sin especificar"
$ws.Range("I177").Value = "Primera"
$ws.Range("J177").Value = 100
$ws.Range("K177").Value = 8000
$ws.Range("L177").Value = 8500
$ws.Range("M177").Value = 8250
$ws.Range("N177").Value = "`$/saco 20 kilos"
$ws.Range("O177").Value = "Provincia de Diguillín"
$ws.Range("P177").Value = 412
$ws.Range("Q177").Value = 20
$ws.Range("R177").Value = "Hortaliza"
